# Sync attendance_reports: rotate "Recorded By" (column G) entries so that
# the last comma-separated contributor is moved to the front of the list.
#
# Example: "dnasr281@gmail.com, System" -> "System, dnasr281@gmail.com"
#          "system, backup@backdoor.com, System" -> "System, system, backup@backdoor.com"
# Single-value cells (no comma) are left unchanged.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$usedRange = $ws.UsedRange
$lastRow = $usedRange.Rows.Count + $usedRange.Row - 1

for ($r = 2; $r -le $lastRow; $r++) {
    $cell = $ws.Cells.Item($r, 7)   # Column G = "Recorded By"
    $value = $cell.Value2

    if ($value -ne $null -and $value -ne "") {
        $parts = $value -split ",\s*"
        if ($parts.Count -gt 1) {
            $last = $parts[$parts.Count - 1]
            $rest = $parts[0..($parts.Count - 2)]
            $newParts = @($last) + $rest
            $newValue = [string]::Join(", ", $newParts)
            $cell.Value2 = $newValue
        }
    }
}
